# Trading update: 2026-02-17 15:36:27
# Append a new OPEN trade (row 51) to both the "All Trades" and the
# "MarketMaking" worksheets with identical data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 51

    $ws.Cells.Item($row, 1).Value = 50

    # Keep the Date column as literal text ("2026-02-17") instead of
    # letting Excel auto-convert it to a date serial number.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"

    $ws.Cells.Item($row, 3).Value = "15:35:59"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.54
    # Exit Price - no value yet (trade still open)
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.5215569553527
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Exit Reason - no value yet (trade still open)
    $ws.Cells.Item($row, 17).Value = 0
}
